$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "244.42" are not
# auto-converted to floating point numbers by the COM layer; this matches
# the original workbook where every Price cell is an inline/shared string.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Simple price / volume(1h) updates (rows 2-18) ---
$ws.Range("D2").Value = "36.612.88"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.957.68"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "244.42"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "58.78"
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.368"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "0.0860"
$ws.Range("E10").Value = "  +9.97%  "
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("D12").Value = "22.59"
$ws.Range("E12").Value = "  +4.64%  "
$ws.Range("D13").Value = "0.832"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "2.246.33"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "13.72"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "5.26"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "1.957.82"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "36.523.08"
$ws.Range("E18").Value = "  +1.41%  "

# --- Row 19/20 swap: Litecoin <-> ShibaInu (ranking order changed) ---
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0878"
$ws.Range("E19").Value = "  +4.11%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "70.15"
$ws.Range("E20").Value = "  -0.80%  "

# --- Simple price / volume(1h) updates (rows 21-33) ---
$ws.Range("D21").Value = "230.49"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").Value = "5.09"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("D26").Value = "9.43"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("D27").Value = "162.52"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").Value = "0.136"
$ws.Range("E28").Value = "  +12.57%  "
$ws.Range("D29").Value = "19.65"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +6.97%  "
$ws.Range("D32").Value = "4.74"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("D33").Value = "0.0642"
$ws.Range("E33").Value = "  +5.62%  "

# --- Row 34/35 swap: InternetComputer(DFINITY) <-> THORChain (ranking order changed) ---
$ws.Range("B34").Value = "THORChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D34").Value = "6.48"
$ws.Range("E34").Value = "  +7.41%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "4.32"
$ws.Range("E35").Value = "  -0.96%  "

# --- Simple price / volume(1h) updates (rows 36-51) ---
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "3.06"
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("D40").Value = "0.1000"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("D42").Value = "2.88"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "16.41"
$ws.Range("E44").Value = "  +4.31%  "
$ws.Range("D45").Value = "1.04"
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").Value = "1.357.80"
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("D47").Value = "88.67"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("D48").Value = "7.27"
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("D49").Value = "2.82"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "46.14"
$ws.Range("E50").Value = "  +5.51%  "
$ws.Range("D51").Value = "2.138.40"
$ws.Range("E51").Value = "  +0.12%  "

# Reset the explicit text-number format back to Normal style so the saved
# cellXfs/style index matches the original (unstyled) Price cells.
$ws.Range("D2:D51").Style = "Normal"
